$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-19 23:48:25"
$ws.Range("G2").Value = "180 cm"
$ws.Range("I2").Value = "4.8 mm"
$ws.Range("E3").Value = "2026-02-19 23:48:27"
$ws.Range("I3").Value = "6.9 mm"
$ws.Range("E4").Value = "2026-02-19 23:48:30"
$ws.Range("H4").Value = "'55%"
$ws.Range("J4").Value = "1010.6 hPa"
$ws.Range("E5").Value = "2026-02-19 23:48:32"
$ws.Range("I5").Value = "8.6 mm"
$ws.Range("E6").Value = "2026-02-19 23:48:35"
$ws.Range("H6").Value = "'74%"
$ws.Range("J6").Value = "1010.8 hPa"
$ws.Range("O6").Value = "10.1 °C"
$ws.Range("E7").Value = "2026-02-19 23:48:37"
$ws.Range("J7").Value = "1011.8 hPa"
$ws.Range("E8").Value = "2026-02-19 23:48:39"
$ws.Range("J8").Value = "1011.5 hPa"
$ws.Range("E9").Value = "2026-02-19 23:48:42"
$ws.Range("E10").Value = "2026-02-19 23:48:44"
$ws.Range("H10").Value = "'73%"
$ws.Range("N10").Value = "2.9 °C 23:05 TU"
$ws.Range("O10").Value = "9.5 °C"
$ws.Range("E11").Value = "2026-02-19 23:48:47"
$ws.Range("O11").Value = "6.2 °C"
$ws.Range("E12").Value = "2026-02-19 23:48:49"
$ws.Range("O12").Value = "10.9 °C"
$ws.Range("E13").Value = "2026-02-19 23:48:51"
$ws.Range("J13").Value = "1012.1 hPa"
$ws.Range("L13").Value = "51.8 km/h - 73º 23:03 TU"
$ws.Range("O13").Value = "4.6 °C"
$ws.Range("E14").Value = "2026-02-19 23:48:54"
$ws.Range("E15").Value = "2026-02-19 23:48:56"
$ws.Range("H15").Value = "'73%"
$ws.Range("O15").Value = "10.0 °C"
$ws.Range("E16").Value = "2026-02-19 23:48:58"
$ws.Range("I16").Value = "12.2 mm"
$ws.Range("E17").Value = "2026-02-19 23:49:01"
$ws.Range("H17").Value = "'77%"
$ws.Range("E18").Value = "2026-02-19 23:49:03"
$ws.Range("H18").Value = "'63%"
$ws.Range("J18").Value = "1010.9 hPa"
$ws.Range("K18").Value = "10.9 MJ/m2"
$ws.Range("N18").Value = "3.4 °C 23:28 TU"
$ws.Range("O18").Value = "10.9 °C"
$ws.Range("E19").Value = "2026-02-19 23:49:05"
$ws.Range("O19").Value = "5.0 °C"
$ws.Range("E20").Value = "2026-02-19 23:49:08"
$ws.Range("E21").Value = "2026-02-19 23:49:10"
$ws.Range("J21").Value = "1012.1 hPa"
$ws.Range("O21").Value = "6.6 °C"
$ws.Range("E22").Value = "2026-02-19 23:49:13"
$ws.Range("I22").Value = "3.2 mm"
$ws.Range("E23").Value = "2026-02-19 23:49:15"
$ws.Range("I23").Value = "12.7 mm"
$ws.Range("E24").Value = "2026-02-19 23:49:17"
$ws.Range("J24").Value = "1015.7 hPa"
$ws.Range("E25").Value = "2026-02-19 23:49:20"
$ws.Range("H25").Value = "'69%"
$ws.Range("I25").Value = "8.2 mm"
$ws.Range("E26").Value = "2026-02-19 23:49:22"
$ws.Range("J26").Value = "1010.5 hPa"
$ws.Range("E27").Value = "2026-02-19 23:49:24"
$ws.Range("E28").Value = "2026-02-19 23:49:27"
$ws.Range("J28").Value = "1010.7 hPa"
$ws.Range("O28").Value = "8.8 °C"
$ws.Range("E29").Value = "2026-02-19 23:49:29"
$ws.Range("K29").Value = "12.4 MJ/m2"
$ws.Range("N29").Value = "3.5 °C 23:24 TU"
$ws.Range("O29").Value = "9.8 °C"
$ws.Range("E30").Value = "2026-02-19 23:49:32"
$ws.Range("J30").Value = "1010.8 hPa"
$ws.Range("O30").Value = "9.8 °C"
$ws.Range("E31").Value = "2026-02-19 23:49:34"
$ws.Range("J31").Value = "1010.0 hPa"
$ws.Range("E32").Value = "2026-02-19 23:49:37"
$ws.Range("H32").Value = "'69%"
$ws.Range("E33").Value = "2026-02-19 23:49:39"
$ws.Range("J33").Value = "1011.7 hPa"
$ws.Range("E34").Value = "2026-02-19 23:49:42"
$ws.Range("E35").Value = "2026-02-19 23:49:44"
$ws.Range("J35").Value = "1017.2 hPa"
$ws.Range("E36").Value = "2026-02-19 23:49:46"
$ws.Range("J36").Value = "1011.0 hPa"
$ws.Range("K36").Value = "10.7 MJ/m2"
$ws.Range("L36").Value = "38.2 km/h - 321º 23:28 TU"
$ws.Range("O36").Value = "11.7 °C"
$ws.Range("E37").Value = "2026-02-19 23:49:48"
$ws.Range("J37").Value = "1012.2 hPa"
$ws.Range("E38").Value = "2026-02-19 23:49:51"
$ws.Range("O38").Value = "11.7 °C"
$ws.Range("E39").Value = "2026-02-19 23:49:53"
$ws.Range("I39").Value = "5.3 mm"
$ws.Range("E40").Value = "2026-02-19 23:49:56"
$ws.Range("H40").Value = "'68%"
$ws.Range("J40").Value = "1013.3 hPa"
$ws.Range("O40").Value = "6.7 °C"
$ws.Range("E41").Value = "2026-02-19 23:49:58"
$ws.Range("J41").Value = "1013.6 hPa"
$ws.Range("E42").Value = "2026-02-19 23:50:00"
$ws.Range("O42").Value = "10.6 °C"
$ws.Range("E43").Value = "2026-02-19 23:50:02"
$ws.Range("N43").Value = "3.4 °C 23:25 TU"
$ws.Range("O43").Value = "8.8 °C"
$ws.Range("E44").Value = "2026-02-19 23:50:05"
$ws.Range("I44").Value = "11.0 mm"
$ws.Range("E45").Value = "2026-02-19 23:50:07"
$ws.Range("I45").Value = "3.9 mm"
$ws.Range("J45").Value = "1016.6 hPa"
$ws.Range("E46").Value = "2026-02-19 23:50:10"
$ws.Range("J46").Value = "1016.5 hPa"
